$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-12 Saturday" "2025-04-13 Sunday"
Replace-Text "423×7=" "231×6="
Replace-Text "101×6=" "392×9="
Replace-Text "113×2=" "102×2="
Replace-Text "312×5=" "573×5="
Replace-Text "720×7=" "762×9="
Replace-Text "839×8=" "342×6="
Replace-Text "374×4=" "759×8="
Replace-Text "854×4=" "490×3="
Replace-Text "607×2=" "547×6="
Replace-Text "267×3=" "368×3="
Replace-Text "827×8=" "816×7="
Replace-Text "882×2=" "450×5="
Replace-Text "600×3=" "999×2="
Replace-Text "551×2=" "917×9="
Replace-Text "638×7=" "512×4="
Replace-Text "561×8=" "268×9="
Replace-Text "999×8=" "266×7="
Replace-Text "925×9=" "612×4="
Replace-Text "243×5=" "433×7="
Replace-Text "613×8=" "616×6="
Replace-Text "395×3=" "808×7="
Replace-Text "307×3=" "765×7="
Replace-Text "392×4=" "530×7="
Replace-Text "663×2=" "177×4="
Replace-Text "631×8=" "431×7="
